$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Header row (bold)
# -----------------------------------------------------------------
$ws.Range("A1").Value = "Comandos mas utilizados"
$ws.Range("B1").Value = "Funcion"
$ws.Range("A1:B1").Font.Bold = $true

# -----------------------------------------------------------------
# Data rows
# -----------------------------------------------------------------
$ws.Range("A2").Value = "git clone"
$ws.Range("B2").Value = "clonar repo"

$ws.Range("A3").Value = "git add . "
$ws.Range("B3").Value = "adherir los archivos que realizamos cambios al stage"

$ws.Range("A4").Value = 'git commit - m "texto commit"'
$ws.Range("B4").Value = "Agregar un commit para pushear"

$ws.Range("A5").Value = "git push origin xxx"
$ws.Range("B5").Value = "Subir los archivos a el repo remoto"

$ws.Range("A6").Value = "git pull "

# B6 is rich text: normal run + underlined run
$b6text = "Actulizar nuestros archivos con el contenido del repo remoto"
$ws.Range("B6").Value = $b6text
$boldEnd = "Actulizar nuestros archivos".Length
$total = $b6text.Length
$ws.Range("B6").Characters($boldEnd + 1, $total - $boldEnd).Font.Underline = $true

# -----------------------------------------------------------------
# Column widths
# -----------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 52.75
$ws.Columns.Item(2).ColumnWidth = 55.5

# -----------------------------------------------------------------
# Selection / active cell
# -----------------------------------------------------------------
[void]$ws.Range("A7").Select()

# -----------------------------------------------------------------
# Page setup
# -----------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
